# Apply the "três tests" edit described in the commit message:
#  - verificação de string  -> update the TypeError message text in F3
#  - verificação de menor que 0 / maior que 10 -> E4 becomes the numeric value -1.11
#    (replacing the old textual "[-1, 11]" placeholder)
#  - widen column F so the longer error text fits
#  - move the active selection to D17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TypeError message (now also mentions strings are not allowed)
$ws.Range("F3").Value = 'TypeError("É necessário que seja um número, strings não inclusas")'

# E4 used to hold the text "[-1, 11]"; it now holds the actual numeric test value
$ws.Range("E4").Value = -1.11

# Column F needs to be wider to fit the new, longer error text
$ws.Columns.Item(6).ColumnWidth = 59.7

# Selection moved from D11 to D17
[void]$ws.Range("D17").Select()
